$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Insert a new column before J -----------------------------------------
# This shifts the existing "setFlags : Bool" column (J) to K and the
# existing "ShiftDirection" column (K) to L, creating a blank J column that
# inherits the centred/wrapped formatting used throughout the header row.
$ws.Columns("J:J").Insert()

# --- New instruction row (row 11): subtractWithCarryS ----------------------
# Populate the text columns first so the new shared-string table entries are
# created in the same order the author's edit produced them (rows before the
# new column header).
$ws.Range("A11").Value = "subtractWithCarryS"
$ws.Range("B11").Value = "Implements various subtraction routines"
$ws.Range("C11").Value = "RegisterID -> RegisterID -> Operand -> MachineState -> bool -> bool -> bool -> ShiftDirection -> MachineState"

$checkmark = [char]0x2713
$ws.Range("D11").Value = $checkmark
$ws.Range("E11").Value = $checkmark
$ws.Range("F11").Value = $checkmark
$ws.Range("H11").Value = $checkmark
$ws.Range("I11").Value = $checkmark
$ws.Range("J11").Value = $checkmark
$ws.Range("K11").Value = $checkmark
$ws.Range("L11").Value = $checkmark

# Match the centred / word-wrapped style used by every other data row (skip
# the empty G11 cell, which this instruction doesn't use).
$dataCells = $ws.Range("A11:F11,H11:L11")
$dataCells.HorizontalAlignment = -4108
$dataCells.WrapText = $true
$ws.Rows("11:11").RowHeight = 64

# --- New column header (J1): "reverse: Bool" --------------------------------
# Set this last so it becomes the final new shared-string entry, matching
# the author's edit ordering.
$ws.Range("J1").Value = "reverse: Bool"

# --- Column width tweaks -----------------------------------------------------
# The new "name" column (A) needs to be a bit wider to fit
# "subtractWithCarryS", and the new "reverse: Bool" column (J) gets its own
# width too; both lose their old bestFit flag once set explicitly.
$ws.Columns("A:A").ColumnWidth = 17.75
$ws.Columns("J:J").ColumnWidth = 11

# --- Selection: land on the newly-added cell --------------------------------
$ws.Range("J11").Select()
